$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the file-id bearing cell values (old guid -> new guid) ---

# Overview sheet
$wsOverview.Range("A2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.md"
$wsOverview.Range("B2").Value = "e2e\17ae3688-f601-4fda-9ec4-75a4b1907617.md"
$wsOverview.Range("G2").Value = "2016-09-01 23:02:56"

# zh-cn sheet
$wsZhCn.Range("A2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.md"
$wsZhCn.Range("G2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.e388192d03ef0a6dd716f17b66d359b14047db4d.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 23:02:51"

# de-de sheet
$wsDeDe.Range("A2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.md"
$wsDeDe.Range("G2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.e388192d03ef0a6dd716f17b66d359b14047db4d.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 23:02:56"

# --- Update the hyperlink display text to match the new file name (in place, via collection iteration) ---

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\17ae3688-f601-4fda-9ec4-75a4b1907617.md"
}

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "17ae3688-f601-4fda-9ec4-75a4b1907617.md"
}

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "17ae3688-f601-4fda-9ec4-75a4b1907617.md"
}
